$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 613, shifting rows 613:654 down to 614:655.
$ws.Rows("613").Insert()

# Populate the new row 613 with the inserted data point (2026/01/11, Sun, 13:00, 201).
# Column A holds a date-looking value that must stay as plain text (matching the
# rest of the sheet, which stores dates as text, not Excel date serials), so we
# temporarily force a Text number format, assign the value, then clear the
# formatting back to the sheet's default (General) so the cell ends up with no
# explicit style -- exactly like every other data row in the sheet.
$ws.Range("A613").NumberFormat = "@"
$ws.Range("A613").Value = "2026/01/11"
$ws.Range("A613").ClearFormats()

$ws.Range("B613").Value = "日"
$ws.Range("C613").Value = 13
$ws.Range("D613").Value = 201
